$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task header: Task 6 : Week 5
$ws.Range("A33").Value = "Task 6 : Week 5"
$ws.Range("A33").Font.Bold = $true

# Fix B27 value from "10 hrs" to "10 hr"
$ws.Range("B27").Value = "10 hr"

$ws.Range("A34").Value = "Keep a log of suggestions"
$ws.Range("B34").Value = "3 hrs"

$ws.Range("A35").Value = "Listen to feedback and make necessary changes"
$ws.Range("B35").Value = "5 hrs"

# New task header: Task 7 : Week 6
$ws.Range("A36").Value = "Task 7 : Week 6"
$ws.Range("A36").Font.Bold = $true

$ws.Range("A37").Value = "Research proxy / alert system"
$ws.Range("A38").Value = "Set up a alert program"
$ws.Range("A39").Value = "Set a timer for alert system"

$ws.Range("B37").Value = "19 hrs"
$ws.Range("B38").Value = "3 hrs"
$ws.Range("B39").Value = "3 hrs"

# New task header: Task 8: Week 7
$ws.Range("A40").Value = "Task 8: Week 7"
$ws.Range("A40").Font.Bold = $true

# Update the view to scroll down and select the newly added last cell
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("A40").Select()
